$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (remove trailing punctuation from header labels)
$ws.Range("D1").Value = "Datafile"
$ws.Range("C1").Value = "wt"
$ws.Range("A1").Value = "Sample"
$ws.Range("B1").Value = "Control"

# Update the selected cell to match the saved selection in the workbook
$ws.Range("D18").Select()
